# Reorder the "COMPETENCES TECHNIQUES" skill lines.
#
# Before (document order):
#   1 Web : client
#   2 Langages : python, matlab, c, c++
#   3 Bases de données : SQL, MongoDB, Neo4j, Redis
#   4 Visualisation : data engineering, tableau
#   5 ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn
#   6 MLOps : node.js, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
#
# After (document order):
#   1 Langages : python, matlab, c, c++
#   2 Visualisation : data engineering, tableau
#   3 MLOps : node.js, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
#   4 Web : client
#   5 ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn
#   6 Bases de données : SQL, MongoDB, Neo4j, Redis
#
# i.e. slot 1 <- old "Langages" text, slot 2 <- old "Visualisation" text,
# slot 3 <- old "MLOps" text, slot 4 <- old "Web" text, slot 5 <- old
# "ML/AI" text (unchanged), slot 6 <- old "Bases de données" text.
#
# This is done in two passes: first every slot's current text is
# replaced by a unique temporary marker (so no slot keeps text that any
# other step might search for), then every marker is replaced by the
# slot's final text.

$d = $word.ActiveDocument

# --- Pass 1: tag each of the 6 slots with a unique temporary marker,
#     based on the slot's CURRENT text (found in document order). ---
$d.Content.Find.Execute("Web : client", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_1__", 2)

$d.Content.Find.Execute("Langages : python, matlab, c, c++", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_2__", 2)

$d.Content.Find.Execute("Bases de données : SQL, MongoDB, Neo4j, Redis", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_3__", 2)

$d.Content.Find.Execute("Visualisation : data engineering, tableau", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_4__", 2)

$d.Content.Find.Execute("ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_5__", 2)

$d.Content.Find.Execute("MLOps : node.js, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__SLOT_6__", 2)

# --- Pass 2: write each slot's new, final content. ---
# slot 1 (was "Web : client") -> "Langages : python, matlab, c, c++"
$d.Content.Find.Execute("__SLOT_1__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Langages : python, matlab, c, c++", 2)

# slot 2 (was "Langages : ...") -> "Visualisation : data engineering, tableau"
$d.Content.Find.Execute("__SLOT_2__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Visualisation : data engineering, tableau", 2)

# slot 3 (was "Bases de données : ...") -> "MLOps : node.js, ..."
$d.Content.Find.Execute("__SLOT_3__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MLOps : node.js, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2)

# slot 4 (was "Visualisation : ...") -> "Web : client"
$d.Content.Find.Execute("__SLOT_4__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Web : client", 2)

# slot 5 (was "ML/AI : ...") -> "ML/AI : ..." (unchanged, restore its own text)
$d.Content.Find.Execute("__SLOT_5__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", 2)

# slot 6 (was "MLOps : ...") -> "Bases de données : SQL, MongoDB, Neo4j, Redis"
$d.Content.Find.Execute("__SLOT_6__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bases de données : SQL, MongoDB, Neo4j, Redis", 2)
